$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.945
$ws.Range("A21").Value = -19.953
$ws.Range("A23").Value = -20.131
$ws.Range("D24").Value = -7.327000000000001
$ws.Range("A25").Value = -21.742
$ws.Range("D28").Value = -8.166
$ws.Range("D36").Value = -7.540999999999999
$ws.Range("D45").Value = -7.425
$ws.Range("D48").Value = -7.540999999999999
$ws.Range("D49").Value = -8.276
$ws.Range("D52").Value = -8.001000000000001
$ws.Range("A53").Value = -22.068
$ws.Range("D53").Value = -8.068999999999999
$ws.Range("D54").Value = -8.105
$ws.Range("A57").Value = -22.253
$ws.Range("A59").Value = -22.461
$ws.Range("A69").Value = -21.606
$ws.Range("D70").Value = -6.853
$ws.Range("A79").Value = -20.866
$ws.Range("A83").Value = -21.976
$ws.Range("D86").Value = -8.251999999999999
$ws.Range("D87").Value = -8.234
$ws.Range("A93").Value = -21.508
$ws.Range("D101").Value = -8.016
